$wb = $excel.ActiveWorkbook

# Sheets that contain the "想去人数" (F column) values that changed
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1439
    $ws.Range("F3").Value = 3030
    $ws.Range("F5").Value = 394
}
